$wb = $excel.ActiveWorkbook

# --- Update "Schedule" sheet (rows 3-4) ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("A3").Value = 46071.0625
$wsSchedule.Range("B3").Value = 46071.25
$wsSchedule.Range("C3").Value = 4.5
$wsSchedule.Range("D3").Value = 17.01
$wsSchedule.Range("E3").Value = 636.442014
$wsSchedule.Range("F3").Value = 37.41575626102293
$wsSchedule.Range("A4").Value = 46071.3125
$wsSchedule.Range("C4").Value = 7.5
$wsSchedule.Range("D4").Value = 28.35
$wsSchedule.Range("E4").Value = 250.087149
$wsSchedule.Range("F4").Value = 8.821416190476191

# --- Update "Detailed" sheet (rows 38-97) ---
$wsDetailed = $wb.Worksheets.Item("Detailed")
$wsDetailed.Range("B38").Value = 84.79000000000001
$wsDetailed.Range("B39").Value = 70.36225
$wsDetailed.Range("B40").Value = 90.82671999999999
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 100.95901
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 105.79
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 103.66848
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 95.72272
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("B45").Value = 84.79000000000001
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("B46").Value = 78
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("B47").Value = 72.75097
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("B48").Value = 66.93682
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("B49").Value = 66.324
$wsDetailed.Range("C49").Value = "historical"
$wsDetailed.Range("B50").Value = 73.19
$wsDetailed.Range("B51").Value = 78.35751
$wsDetailed.Range("B52").Value = 78.95034
$wsDetailed.Range("E52").Value = "OFF"
$wsDetailed.Range("B53").Value = 74.77921000000001
$wsDetailed.Range("B54").Value = 73.20017
$wsDetailed.Range("B55").Value = 73.20016
$wsDetailed.Range("B56").Value = 73.2
$wsDetailed.Range("B57").Value = 73.2
$wsDetailed.Range("B58").Value = 73.20017
$wsDetailed.Range("B59").Value = 65
$wsDetailed.Range("B60").Value = 73.20013
$wsDetailed.Range("E60").Value = "ON"
$wsDetailed.Range("B61").Value = 73.7812
$wsDetailed.Range("E61").Value = "ON"
$wsDetailed.Range("B62").Value = 84.79000000000001
$wsDetailed.Range("B63").Value = 110.45944
$wsDetailed.Range("B64").Value = 84.79000000000001
$wsDetailed.Range("E64").Value = "OFF"
$wsDetailed.Range("B65").Value = 35.88
$wsDetailed.Range("B66").Value = 13.52897
$wsDetailed.Range("B69").Value = 0.51
$wsDetailed.Range("B71").Value = 0.51
$wsDetailed.Range("B73").Value = 0.7
$wsDetailed.Range("B74").Value = 22.07
$wsDetailed.Range("B75").Value = 22.07
$wsDetailed.Range("B76").Value = 22.07
$wsDetailed.Range("B77").Value = 36.06046
$wsDetailed.Range("B78").Value = 36.06021
$wsDetailed.Range("B79").Value = 65
$wsDetailed.Range("B80").Value = 79.95016
$wsDetailed.Range("B81").Value = 71.85057999999999
$wsDetailed.Range("B82").Value = 71.25176
$wsDetailed.Range("B83").Value = 66.72439
$wsDetailed.Range("B84").Value = 81.14199000000001
$wsDetailed.Range("B85").Value = 49.38732
$wsDetailed.Range("B86").Value = 56.72622
$wsDetailed.Range("B87").Value = 100.81666
$wsDetailed.Range("B88").Value = 159.33954
$wsDetailed.Range("B89").Value = 151.65295
$wsDetailed.Range("B90").Value = 144.56143
$wsDetailed.Range("B91").Value = 143.3261
$wsDetailed.Range("B92").Value = 138.42
$wsDetailed.Range("B93").Value = 108.89
$wsDetailed.Range("B94").Value = 108.89
$wsDetailed.Range("B95").Value = 105.79
$wsDetailed.Range("B96").Value = 105.79
$wsDetailed.Range("B97").Value = 97.55358
